$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 162 (shifts existing rows 162-266 down to 163-267).
$ws.Rows.Item(162).Insert()

# Populate the newly inserted row 162 with a new data record (same structure as the
# surrounding rows: Mercado/Region/Categoria static columns stay constant).
$ws.Cells.Item(162, 1).Value2 = 5
$ws.Cells.Item(162, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(162, 3).Value2 = "Maule"

$ws.Cells.Item(162, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(162, 4).Value2 = 44582

$ws.Cells.Item(162, 5).Value2 = 7
$ws.Cells.Item(162, 6).Value2 = 100114013
$ws.Cells.Item(162, 7).Value2 = "Zanahoria"
$ws.Cells.Item(162, 8).Value2 = "Sin especificar"
$ws.Cells.Item(162, 9).Value2 = "Primera"
$ws.Cells.Item(162, 10).Value2 = 250
$ws.Cells.Item(162, 11).Value2 = 8000
$ws.Cells.Item(162, 12).Value2 = 8000
$ws.Cells.Item(162, 13).Value2 = 8000
$ws.Cells.Item(162, 14).Value2 = "`$/saco 20 kilos"
$ws.Cells.Item(162, 15).Value2 = "Región de Ñuble"
$ws.Cells.Item(162, 16).Value2 = 400
$ws.Cells.Item(162, 17).Value2 = 20
$ws.Cells.Item(162, 18).Value2 = "Hortaliza"
